$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 879.9
$ws.Range("I62").Value = 866.55554
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 866.55554
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -242.55554
$ws.Range("N62").Value = -2248
$ws.Range("H65").Value = 879.9
$ws.Range("I65").Value = 866.55554
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 4332.7777
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -1212.7777
$ws.Range("N65").Value = -11240
$ws.Range("H116").Value = 2593
$ws.Range("J116").Value = 3170
$ws.Range("L116").Value = 3170
$ws.Range("N116").Value = -10054
$ws.Range("H132").Value = 5955099
$ws.Range("I132").Value = 7521651.5
$ws.Range("J132").Value = 2199.6
$ws.Range("K132").Value = 22564954.5
$ws.Range("L132").Value = 6598.799999999999
$ws.Range("M132").Value = -22562424.5
$ws.Range("N132").Value = -11658.8
$ws.Range("H137").Value = 1302.0702
$ws.Range("J137").Value = 2186.8096
$ws.Range("L137").Value = 6560.4288
$ws.Range("N137").Value = -11660.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16425.44
$ws.Range("I32").Value = 17269.135
$ws.Range("J32").Value = 12752.883
$ws.Range("K32").Value = 17269.135
$ws.Range("L32").Value = 12752.883
$ws.Range("M32").Value = -16982.135
$ws.Range("N32").Value = -13326.883
$ws.Range("H45").Value = 37038748
$ws.Range("I45").Value = 83334690
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 83334690
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -83334313
$ws.Range("N45").Value = -2754
$ws.Range("H61").Value = 1499.4419
$ws.Range("I61").Value = 1163.5807
$ws.Range("J61").Value = 2367.0833
$ws.Range("K61").Value = 1163.5807
$ws.Range("L61").Value = 2367.0833
$ws.Range("M61").Value = -951.5807
$ws.Range("N61").Value = -2791.0833
$ws.Range("H63").Value = 3335100
$ws.Range("I63").Value = 5001250
$ws.Range("J63").Value = 2800
$ws.Range("K63").Value = 5001250
$ws.Range("L63").Value = 2800
$ws.Range("M63").Value = -5000564
$ws.Range("N63").Value = -4172
$ws.Range("H66").Value = 3335100
$ws.Range("I66").Value = 5001250
$ws.Range("J66").Value = 2800
$ws.Range("K66").Value = 25006250
$ws.Range("L66").Value = 14000
$ws.Range("M66").Value = -25002818
$ws.Range("N66").Value = -20864
$ws.Range("H97").Value = 975.86206
$ws.Range("I97").Value = 717.2632
$ws.Range("J97").Value = 1467.2
$ws.Range("K97").Value = 717.2632
$ws.Range("L97").Value = 1467.2
$ws.Range("M97").Value = -221.2632
$ws.Range("N97").Value = -2459.2
$ws.Range("H102").Value = 1494
$ws.Range("I102").Value = 988
$ws.Range("K102").Value = 988
$ws.Range("M102").Value = 634
$ws.Range("H110").Value = 1837
$ws.Range("I110").Value = 2008.75
$ws.Range("J110").Value = 1150
$ws.Range("K110").Value = 2008.75
$ws.Range("L110").Value = 1150
$ws.Range("M110").Value = 36.25
$ws.Range("N110").Value = -5240
$ws.Range("H132").Value = 4277.32
$ws.Range("I132").Value = 4664.243
$ws.Range("J132").Value = 3176.077
$ws.Range("K132").Value = 13992.729
$ws.Range("L132").Value = 9528.231
$ws.Range("M132").Value = -11462.729
$ws.Range("N132").Value = -14588.231
$ws.Range("H136").Value = 1499.4419
$ws.Range("I136").Value = 1163.5807
$ws.Range("J136").Value = 2367.0833
$ws.Range("K136").Value = 3490.7421
$ws.Range("L136").Value = 7101.249899999999
$ws.Range("M136").Value = -940.7420999999999
$ws.Range("N136").Value = -12201.2499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 314
$ws.Range("I80").Value = 297.33334
$ws.Range("J80").Value = 322.33334
$ws.Range("K80").Value = 297.33334
$ws.Range("L80").Value = 322.33334
$ws.Range("M80").Value = 700.66666
$ws.Range("N80").Value = -2318.33334
$ws.Range("H83").Value = 314
$ws.Range("I83").Value = 297.33334
$ws.Range("J83").Value = 322.33334
$ws.Range("K83").Value = 1486.6667
$ws.Range("L83").Value = 1611.6667
$ws.Range("M83").Value = 3505.3333
$ws.Range("N83").Value = -11595.6667
$ws.Range("H86").Value = 2850
$ws.Range("I86").Value = 2760
$ws.Range("K86").Value = 2760
$ws.Range("M86").Value = -1637
$ws.Range("H89").Value = 2850
$ws.Range("I89").Value = 2760
$ws.Range("K89").Value = 13800
$ws.Range("M89").Value = -8184
$ws.Range("H99").Value = 2833.3333
$ws.Range("I99").Value = 2750
$ws.Range("K99").Value = 2750
$ws.Range("M99").Value = -1252
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H107").Value = 1427.381
$ws.Range("I107").Value = 1459
$ws.Range("K107").Value = 1459
$ws.Range("M107").Value = 461
$ws.Range("H134").Value = 1736.8813
$ws.Range("I134").Value = 1509.025
$ws.Range("K134").Value = 4527.075000000001
$ws.Range("M134").Value = -1992.075000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 359.7619
$ws.Range("I22").Value = 351.46155
$ws.Range("J22").Value = 373.25
$ws.Range("K22").Value = 351.46155
$ws.Range("L22").Value = 373.25
$ws.Range("M22").Value = -1.461549999999988
$ws.Range("N22").Value = -1073.25
$ws.Range("H31").Value = 3971180
$ws.Range("I31").Value = 2312.3635
$ws.Range("K31").Value = 2312.3635
$ws.Range("M31").Value = -2017.3635
$ws.Range("H34").Value = 3971180
$ws.Range("I34").Value = 2312.3635
$ws.Range("K34").Value = 2312.3635
$ws.Range("M34").Value = -2110.3635
$ws.Range("H58").Value = 1591.8334
$ws.Range("I58").Value = 1425.25
$ws.Range("J58").Value = 1925
$ws.Range("K58").Value = 1425.25
$ws.Range("L58").Value = 1925
$ws.Range("M58").Value = -1222.25
$ws.Range("N58").Value = -2331
$ws.Range("H99").Value = 1986.1724
$ws.Range("I99").Value = 1699.5
$ws.Range("J99").Value = 2339
$ws.Range("K99").Value = 1699.5
$ws.Range("L99").Value = 2339
$ws.Range("M99").Value = -201.5
$ws.Range("N99").Value = -5335
$ws.Range("H105").Value = 649.2
$ws.Range("I105").Value = 480.77777
$ws.Range("J105").Value = 901.8333
$ws.Range("K105").Value = 480.77777
$ws.Range("L105").Value = 901.8333
$ws.Range("M105").Value = 1266.22223
$ws.Range("N105").Value = -4395.8333
$ws.Range("H107").Value = 575.0909
$ws.Range("I107").Value = 504.78378
$ws.Range("K107").Value = 504.78378
$ws.Range("M107").Value = 1415.21622
$ws.Range("H126").Value = 1986.1724
$ws.Range("I126").Value = 1699.5
$ws.Range("J126").Value = 2339
$ws.Range("K126").Value = 5098.5
$ws.Range("L126").Value = 7017
$ws.Range("M126").Value = -2628.5
$ws.Range("N126").Value = -11957
$ws.Range("H132").Value = 3380771
$ws.Range("I132").Value = 1789.6154
$ws.Range("J132").Value = 11367454
$ws.Range("K132").Value = 5368.8462
$ws.Range("L132").Value = 34102362
$ws.Range("M132").Value = -2838.8462
$ws.Range("N132").Value = -34107422
$ws.Range("H134").Value = 1199.762
$ws.Range("I134").Value = 1168.1578
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 3504.4734
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -969.4733999999999
$ws.Range("N134").Value = -9570
$ws.Range("H136").Value = 1591.8334
$ws.Range("I136").Value = 1425.25
$ws.Range("J136").Value = 1925
$ws.Range("K136").Value = 4275.75
$ws.Range("L136").Value = 5775
$ws.Range("M136").Value = -1725.75
$ws.Range("N136").Value = -10875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 147
$ws.Range("I8").Value = 147
$ws.Range("K8").Value = 441
$ws.Range("M8").Value = -302
$ws.Range("H132").Value = 1434.6786
$ws.Range("I132").Value = 645.8421
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 5812.5789
$ws.Range("L132").Value = 27900
$ws.Range("M132").Value = -3282.5789
$ws.Range("N132").Value = -32960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 751.3214
$ws.Range("I97").Value = 689.875
$ws.Range("K97").Value = 689.875
$ws.Range("M97").Value = -193.875
$ws.Range("H102").Value = 1859.4286
$ws.Range("I102").Value = 2004
$ws.Range("K102").Value = 2004
$ws.Range("M102").Value = -382
$ws.Range("H122").Value = 31252226
$ws.Range("I122").Value = 58825468
$ws.Range("J122").Value = 2553.8667
$ws.Range("K122").Value = 176476404
$ws.Range("L122").Value = 7661.6001
$ws.Range("M122").Value = -176473954
$ws.Range("N122").Value = -12561.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27779658
$ws.Range("I7").Value = 1483.4
$ws.Range("J7").Value = 62502376
$ws.Range("K7").Value = 1483.4
$ws.Range("L7").Value = 62502376
$ws.Range("M7").Value = -1371.4
$ws.Range("N7").Value = -62502600
$ws.Range("H68").Value = 15626581
$ws.Range("I68").Value = 1152.3334
$ws.Range("J68").Value = 25001838
$ws.Range("K68").Value = 1152.3334
$ws.Range("L68").Value = 25001838
$ws.Range("M68").Value = -403.3334
$ws.Range("N68").Value = -25003336
$ws.Range("H71").Value = 15626581
$ws.Range("I71").Value = 1152.3334
$ws.Range("J71").Value = 25001838
$ws.Range("K71").Value = 5761.666999999999
$ws.Range("L71").Value = 125009190
$ws.Range("M71").Value = -2017.666999999999
$ws.Range("N71").Value = -125016678
$ws.Range("H126").Value = 27779658
$ws.Range("I126").Value = 1483.4
$ws.Range("J126").Value = 62502376
$ws.Range("K126").Value = 4450.200000000001
$ws.Range("L126").Value = 187507128
$ws.Range("M126").Value = -1980.200000000001
$ws.Range("N126").Value = -187512068
$ws.Range("H132").Value = 2518.4285
$ws.Range("I132").Value = 2232.1538
$ws.Range("J132").Value = 2766.5334
$ws.Range("K132").Value = 6696.4614
$ws.Range("L132").Value = 8299.600199999999
$ws.Range("M132").Value = -4166.4614
$ws.Range("N132").Value = -13359.6002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 62500400
$ws.Range("I96").Value = 62500400
$ws.Range("K96").Value = 62500400
$ws.Range("M96").Value = -62499027
$ws.Range("H136").Value = 1095.0193
$ws.Range("I136").Value = 1111
$ws.Range("J136").Value = 1059.0625
$ws.Range("K136").Value = 3333
$ws.Range("L136").Value = 3177.1875
$ws.Range("M136").Value = -783
$ws.Range("N136").Value = -8277.1875
